$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Copy the "closing" (bottom-border) row format from the last
#    worker row (21, CARLOS - about to be removed) onto row 20
#    (LINA), which will become the new last worker row once the
#    rows below it are deleted.
# ------------------------------------------------------------------
$ws.Range("B21:J21").Copy() | Out-Null
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Remove the workers that are no longer part of the updated
#    database: ANGELYS DAYENNER MATOS PACHECO (row 16),
#    ANGIE LUZ MIRANDA BELLO (row 19) and CARLOS ARTURO ARRIETA
#    SAENZ (row 21). Deleting from the bottom up keeps the row
#    numbers of the rows still to be removed stable.
# ------------------------------------------------------------------
$ws.Rows("21").EntireRow.Delete() | Out-Null
$ws.Rows("19").EntireRow.Delete() | Out-Null
$ws.Rows("16").EntireRow.Delete() | Out-Null

# ------------------------------------------------------------------
# 3. Update the remaining 3 workers' "Periodo Mora" (2507 -> 2508)
#    and refresh the "Salario Basico" amounts for the two workers
#    whose value changed.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"

$ws.Range("G17").Value = 1962956
$ws.Range("G18").Value = 1656244

# ------------------------------------------------------------------
# 4. Update the summary figures: total "Valor Mora" and the
#    worker head-count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 201708
$ws.Range("C13").Value = 3

# ------------------------------------------------------------------
# 5. The "Nombre Trabajador" column shrinks now that the longest
#    name in the sheet is shorter (best-fit/autofit effect).
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 29.33
